$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hide the columns that carry the intermediate calculations (B:I, J, K:N, V)
$ws.Range("B1:I1").EntireColumn.Hidden = $true
$ws.Range("J1:J1").EntireColumn.Hidden = $true
$ws.Range("K1:N1").EntireColumn.Hidden = $true
$ws.Range("V1:V1").EntireColumn.Hidden = $true

# Move the active selection to U22
$ws.Range("U22").Select() | Out-Null
